$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- LASER CANNON (row 9): Reload "13/12/11" -> "10-10-10", Reload M 12 -> 10 ---
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "10-10-10"
$ws.Range("E9").NumberFormat = "YYYY/MM/DD"
$ws.Range("F9").Value = 10

# --- PLASMA BEAM (row 10): Reload "16/15/14" -> "11/11/11", Reload M 15 -> 11 ---
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11/11/11"
$ws.Range("E10").NumberFormat = "YYYY/MM/DD"
$ws.Range("F10").Value = 11

# --- Update the view state left from the user's last selection ---
$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollColumn = 3
